$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01900254227964501
$ws.Range("C2").Value = 0.2096779636302306

$ws.Range("B3").Value = 0.05681013924720223
$ws.Range("C3").Value = 0.2182627364139333

$ws.Range("B4").Value = 0.8766164083541255
$ws.Range("C4").Value = 0.4433754341917704

$ws.Range("B5").Value = 0.9933820554737847
$ws.Range("C5").Value = 0.3818292616081996

$ws.Range("B6").Value = 0.9802102477359591
$ws.Range("C6").Value = 0.8092171988767494

$ws.Range("B7").Value = 0.9528082279664156
$ws.Range("C7").Value = 0.3023678360306385

$ws.Range("B8").Value = 0.009752786755561828
$ws.Range("C8").Value = 0.1822593879699707

$ws.Range("B9").Value = 0.1678093798411143
$ws.Range("C9").Value = 0.2054546806956996

$ws.Range("B10").Value = 0.6938283950552282
$ws.Range("C10").Value = 0.4134112841906342
